$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44830
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 30000
$ws.Range("O2").Value = 30000
$ws.Range("P2").Value = 30000
$ws.Range("R2").Value = "Perú"
$ws.Range("S2").Value = 1500

$ws.Range("D3").Value = 44333
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 22000
$ws.Range("R3").Value = "Perú"
$ws.Range("S3").Value = 1100

$ws.Range("D4").Value = 44438
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 21000
$ws.Range("R4").Value = "Perú"
$ws.Range("S4").Value = 1050

$ws.Range("D5").Value = 44613
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 30000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 30000
$ws.Range("R5").Value = "Perú"
$ws.Range("S5").Value = 1500

$ws.Range("D6").Value = 44445
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 35
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("R6").Value = "Perú"
$ws.Range("S6").Value = 1000

$ws.Range("D7").Value = 44277
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 24000
$ws.Range("P7").Value = 24000
$ws.Range("R7").Value = "Perú"
$ws.Range("S7").Value = 1200

$ws.Range("D8").Value = 44529
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 34
$ws.Range("N8").Value = 28000
$ws.Range("O8").Value = 28000
$ws.Range("P8").Value = 28000
$ws.Range("R8").Value = "Perú"
$ws.Range("S8").Value = 1400

$ws.Range("D9").Value = 44620
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 22000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 22000
$ws.Range("R9").Value = "Perú"
$ws.Range("S9").Value = 1100

$ws.Range("D10").Value = 44350
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 90
$ws.Range("N10").Value = 21000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 21556
$ws.Range("R10").Value = "Perú"
$ws.Range("S10").Value = 1078

$ws.Range("D11").Value = 44442
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 22000
$ws.Range("O11").Value = 22000
$ws.Range("P11").Value = 22000
$ws.Range("R11").Value = "Perú"
$ws.Range("S11").Value = 1100

$ws.Range("D12").Value = 44357
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 21000
$ws.Range("P12").Value = 20500
$ws.Range("R12").Value = "Perú"
$ws.Range("S12").Value = 1025

$ws.Range("D13").Value = 44302
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 19000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 19500
$ws.Range("R13").Value = "Perú"
$ws.Range("S13").Value = 975

$ws.Range("D14").Value = 44452
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 35
$ws.Range("N14").Value = 21000
$ws.Range("O14").Value = 22000
$ws.Range("P14").Value = 21429
$ws.Range("R14").Value = "Perú"
$ws.Range("S14").Value = 1071

$ws.Range("D15").Value = 44165
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 27000
$ws.Range("O15").Value = 28000
$ws.Range("P15").Value = 27500
$ws.Range("R15").Value = "Perú"
$ws.Range("S15").Value = 1375

$ws.Range("D16").Value = 44354
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = 21000
$ws.Range("O16").Value = 22000
$ws.Range("P16").Value = 21500
$ws.Range("R16").Value = "Perú"
$ws.Range("S16").Value = 1075

$ws.Range("D17").Value = 44396
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 45
$ws.Range("N17").Value = 22000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 22000
$ws.Range("R17").Value = "Perú"
$ws.Range("S17").Value = 1100

$ws.Range("D18").Value = 44473
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 40
$ws.Range("N18").Value = 24000
$ws.Range("O18").Value = 24000
$ws.Range("P18").Value = 24000
$ws.Range("R18").Value = "Perú"
$ws.Range("S18").Value = 1200

$ws.Range("D19").Value = 44312
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 22000
$ws.Range("O19").Value = 22000
$ws.Range("P19").Value = 22000
$ws.Range("R19").Value = "Perú"
$ws.Range("S19").Value = 1100

$ws.Range("D20").Value = 44363
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 150
$ws.Range("N20").Value = 21000
$ws.Range("O20").Value = 22000
$ws.Range("P20").Value = 21500
$ws.Range("R20").Value = "Perú"
$ws.Range("S20").Value = 1075

$ws.Range("D21").Value = 44270
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 24000
$ws.Range("O21").Value = 24000
$ws.Range("P21").Value = 24000
$ws.Range("R21").Value = "Perú"
$ws.Range("S21").Value = 1200

$ws.Range("D22").Value = 44166
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 120
$ws.Range("N22").Value = 28000
$ws.Range("O22").Value = 28000
$ws.Range("P22").Value = 28000
$ws.Range("R22").Value = "Perú"
$ws.Range("S22").Value = 1400

$ws.Range("D23").Value = 44435
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 60
$ws.Range("N23").Value = 25000
$ws.Range("O23").Value = 25000
$ws.Range("P23").Value = 25000
$ws.Range("R23").Value = "Perú"
$ws.Range("S23").Value = 1250

$ws.Range("D24").Value = 44326
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 40
$ws.Range("N24").Value = 22000
$ws.Range("O24").Value = 22000
$ws.Range("P24").Value = 22000
$ws.Range("R24").Value = "Perú"
$ws.Range("S24").Value = 1100

$ws.Range("D25").Value = 44424
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 70
$ws.Range("N25").Value = 24000
$ws.Range("O25").Value = 25000
$ws.Range("P25").Value = 24429
$ws.Range("R25").Value = "Perú"
$ws.Range("S25").Value = 1221

$ws.Range("D26").Value = 44263
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 150
$ws.Range("N26").Value = 15000
$ws.Range("O26").Value = 15000
$ws.Range("P26").Value = 15000
$ws.Range("R26").Value = "Perú"
$ws.Range("S26").Value = 750

$ws.Range("D27").Value = 44299
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 150
$ws.Range("N27").Value = 19000
$ws.Range("O27").Value = 20000
$ws.Range("P27").Value = 19500
$ws.Range("R27").Value = "Perú"
$ws.Range("S27").Value = 975

$ws.Range("D28").Value = 44305
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 40
$ws.Range("N28").Value = 24000
$ws.Range("O28").Value = 24000
$ws.Range("P28").Value = 24000
$ws.Range("R28").Value = "Perú"
$ws.Range("S28").Value = 1200

$ws.Range("D29").Value = 44300
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 150
$ws.Range("N29").Value = 19000
$ws.Range("O29").Value = 20000
$ws.Range("P29").Value = 19500
$ws.Range("R29").Value = "Perú"
$ws.Range("S29").Value = 975

$ws.Range("D30").Value = 44417
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 30
$ws.Range("N30").Value = 24000
$ws.Range("O30").Value = 24000
$ws.Range("P30").Value = 24000
$ws.Range("R30").Value = "Perú"
$ws.Range("S30").Value = 1200

$ws.Range("D31").Value = 44760
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 300
$ws.Range("N31").Value = 24000
$ws.Range("O31").Value = 25000
$ws.Range("P31").Value = 24500
$ws.Range("R31").Value = "Perú"
$ws.Range("S31").Value = 1225

$ws.Range("D32").Value = 44365
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 150
$ws.Range("N32").Value = 20000
$ws.Range("O32").Value = 21000
$ws.Range("P32").Value = 20500
$ws.Range("R32").Value = "Perú"
$ws.Range("S32").Value = 1025

$ws.Range("D33").Value = 44382
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 200
$ws.Range("N33").Value = 19000
$ws.Range("O33").Value = 20000
$ws.Range("P33").Value = 19500
$ws.Range("R33").Value = "Perú"
$ws.Range("S33").Value = 975

$ws.Range("D34").Value = 44356
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 100
$ws.Range("N34").Value = 20000
$ws.Range("O34").Value = 21000
$ws.Range("P34").Value = 20500
$ws.Range("R34").Value = "Perú"
$ws.Range("S34").Value = 1025

$ws.Range("D35").Value = 44372
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 60
$ws.Range("N35").Value = 20000
$ws.Range("O35").Value = 21000
$ws.Range("P35").Value = 20667
$ws.Range("R35").Value = "Perú"
$ws.Range("S35").Value = 1033

$ws.Range("D36").Value = 44298
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 240
$ws.Range("N36").Value = 19000
$ws.Range("O36").Value = 20000
$ws.Range("P36").Value = 19500
$ws.Range("R36").Value = "Perú"
$ws.Range("S36").Value = 975

$ws.Range("D37").Value = 44522
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 25
$ws.Range("N37").Value = 30000
$ws.Range("O37").Value = 30000
$ws.Range("P37").Value = 30000
$ws.Range("R37").Value = "Perú"
$ws.Range("S37").Value = 1500

$ws.Range("D38").Value = 44355
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 200
$ws.Range("N38").Value = 20000
$ws.Range("O38").Value = 21000
$ws.Range("P38").Value = 20500
$ws.Range("R38").Value = "Ecuador"
$ws.Range("S38").Value = 1025

$ws.Range("D39").Value = 44410
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 40
$ws.Range("N39").Value = 25000
$ws.Range("O39").Value = 25000
$ws.Range("P39").Value = 25000
$ws.Range("R39").Value = "Perú"
$ws.Range("S39").Value = 1250

$ws.Range("D40").Value = 44431
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 60
$ws.Range("N40").Value = 25000
$ws.Range("O40").Value = 25000
$ws.Range("P40").Value = 25000
$ws.Range("R40").Value = "Perú"
$ws.Range("S40").Value = 1250

$ws.Range("D41").Value = 44284
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 40
$ws.Range("N41").Value = 23000
$ws.Range("O41").Value = 23000
$ws.Range("P41").Value = 23000
$ws.Range("R41").Value = "Perú"
$ws.Range("S41").Value = 1150
